$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.676499999999995
$ws.Range("B4").Value = 4.887200000000004
$ws.Range("C4").Value = -14.2594
$ws.Range("D4").Value = -7.735499999999996
$ws.Range("C5").Value = -14.49870000000002
$ws.Range("B6").Value = 8.979400000000007
$ws.Range("C6").Value = -12.05830000000001
$ws.Range("B7").Value = 5.086399999999997
$ws.Range("B8").Value = 4.771
$ws.Range("C8").Value = -11.38259999999999
$ws.Range("D9").Value = -8.137600000000003
$ws.Range("D11").Value = -8.372199999999998
$ws.Range("D14").Value = -7.223299999999996
$ws.Range("B16").Value = 8.871200000000007
$ws.Range("C16").Value = -12.0622
$ws.Range("D18").Value = -8.268499999999994
$ws.Range("B20").Value = 5.136900000000002
$ws.Range("B21").Value = 5.037099999999996
$ws.Range("C22").Value = -10.88069999999999
$ws.Range("D25").Value = -8.456599999999996
